$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B51").Value = 3.5
$ws.Range("C51").Value = 3

$ws.Range("B52").Value = 3.4

$ws.Range("B53").Value = 3.2

$ws.Range("B56").Value = 3.2

$ws.Range("C59").Value = 3.9

$ws.Range("C60").Value = 3.6

$ws.Range("C63").Value = 4.3

$ws.Range("C64").Value = 0.8

$ws.Range("C65").Value = 2.2

$ws.Range("B66").Value = 3.1
$ws.Range("C66").Value = 5.6

$ws.Range("B67").Value = 3.5
$ws.Range("C67").Value = 5.1

$ws.Range("B68").Value = 2.7
$ws.Range("C68").Value = 0.2

$ws.Range("C69").Value = 0.1

$ws.Range("B70").Value = -2.8
$ws.Range("C70").Value = -16.3

$ws.Range("B71").Value = -8.9
$ws.Range("C71").Value = -8.5

$ws.Range("C72").Value = 22.9

$ws.Range("B73").Value = -2.2
$ws.Range("C73").Value = 12.8
